$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: quarter period labels (shifted left, new quarter appended)
$ws.Cells.Item(8, 4).Value = "فصل سوم منتهی به 1399/09"
$ws.Cells.Item(8, 5).Value = "فصل چهارم منتهی به 1399/12"
$ws.Cells.Item(8, 6).Value = "فصل اول منتهی به 1400/03"
$ws.Cells.Item(8, 7).Value = "فصل دوم منتهی به 1400/06"
$ws.Cells.Item(8, 8).Value = "فصل سوم منتهی به 1400/09"
$ws.Cells.Item(8, 9).Value = "فصل چهارم منتهی به 1400/12"
$ws.Cells.Item(8, 10).Value = "فصل اول منتهی به 1401/03"
$ws.Cells.Item(8, 11).Value = "فصل دوم منتهی به 1401/06"
$ws.Cells.Item(8, 12).Value = "فصل سوم منتهی به 1401/09"
$ws.Cells.Item(8, 13).Value = "فصل چهارم منتهی به 1401/12"

# Row 9: publish dates (shifted left, new date appended)
$ws.Cells.Item(9, 4).Value = "1400-12-23 (3)"
$ws.Cells.Item(9, 5).Value = "1401-04-05 (10)"
$ws.Cells.Item(9, 6).Value = "1401-04-29 (3)"
$ws.Cells.Item(9, 7).Value = "1401-08-30 (4)"
$ws.Cells.Item(9, 8).Value = "1401-10-29 (3)"
$ws.Cells.Item(9, 9).Value = "1402-02-27 (7)"
$ws.Cells.Item(9, 10).Value = "1401-04-29"
$ws.Cells.Item(9, 11).Value = "1401-08-30 (2)"
$ws.Cells.Item(9, 12).Value = "1401-10-29"
$ws.Cells.Item(9, 13).Value = "1402-02-27"

# Row 12
$ws.Cells.Item(12, 4).Value = -122929
$ws.Cells.Item(12, 5).Value = -42747
$ws.Cells.Item(12, 6).Value = 480890
$ws.Cells.Item(12, 7).Value = -557697
$ws.Cells.Item(12, 8).Value = -420719
$ws.Cells.Item(12, 9).Value = -365781
$ws.Cells.Item(12, 10).Value = -794141
$ws.Cells.Item(12, 11).Value = -937143
$ws.Cells.Item(12, 12).Value = 1127987
$ws.Cells.Item(12, 13).Value = 6062

# Row 13
$ws.Cells.Item(13, 4).Value = -2000
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = -3445
$ws.Cells.Item(13, 8).Value = -24127
$ws.Cells.Item(13, 9).Value = -18210
$ws.Cells.Item(13, 10).Value = -9875
$ws.Cells.Item(13, 11).Value = -55936
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = -29370

# Row 14
$ws.Cells.Item(14, 4).Value = -124929
$ws.Cells.Item(14, 5).Value = -42747
$ws.Cells.Item(14, 6).Value = 480890
$ws.Cells.Item(14, 7).Value = -561142
$ws.Cells.Item(14, 8).Value = -444846
$ws.Cells.Item(14, 9).Value = -383991
$ws.Cells.Item(14, 10).Value = -804016
$ws.Cells.Item(14, 11).Value = -993079
$ws.Cells.Item(14, 12).Value = 1127987
$ws.Cells.Item(14, 13).Value = -23308

# Row 16
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = 0

# Row 17
$ws.Cells.Item(17, 4).Value = -522
$ws.Cells.Item(17, 5).Value = -45136
$ws.Cells.Item(17, 6).Value = -436011
$ws.Cells.Item(17, 7).Value = 344715
$ws.Cells.Item(17, 8).Value = -22378
$ws.Cells.Item(17, 9).Value = -248876
$ws.Cells.Item(17, 10).Value = -1653
$ws.Cells.Item(17, 11).Value = -106173
$ws.Cells.Item(17, 12).Value = -198457
$ws.Cells.Item(17, 13).Value = -1307409

# Row 18
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = 0

# Row 19
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = 0

# Row 20
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = 0

# Row 21
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 13).Value = 0

# Row 22
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 0
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = 0

# Row 23
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 700
$ws.Cells.Item(23, 12).Value = 1250
$ws.Cells.Item(23, 13).Value = -35

# Row 24
$ws.Cells.Item(24, 4).Value = -11500
$ws.Cells.Item(24, 5).Value = 9497
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = -549
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 9).Value = 500
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = -49
$ws.Cells.Item(24, 12).Value = 49
$ws.Cells.Item(24, 13).Value = 0

# Row 25
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 302750
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 25450
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 2601

# Row 26
$ws.Cells.Item(26, 4).Value = 11500
$ws.Cells.Item(26, 5).Value = -191071
$ws.Cells.Item(26, 6).Value = -6700
$ws.Cells.Item(26, 7).Value = -846830
$ws.Cells.Item(26, 8).Value = 264000
$ws.Cells.Item(26, 9).Value = 589530
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 0
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = 0

# Row 27
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 8).Value = 0
$ws.Cells.Item(27, 9).Value = 0
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = 0

# Row 28
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0

# Row 29
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(29, 9).Value = 0
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 13).Value = 0

# Row 30
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0

# Row 31
$ws.Cells.Item(31, 4).Value = 9142
$ws.Cells.Item(31, 5).Value = 7696
$ws.Cells.Item(31, 6).Value = 15598
$ws.Cells.Item(31, 7).Value = 31477
$ws.Cells.Item(31, 8).Value = 54288
$ws.Cells.Item(31, 9).Value = 29989
$ws.Cells.Item(31, 10).Value = 2021
$ws.Cells.Item(31, 11).Value = 1030
$ws.Cells.Item(31, 12).Value = 676
$ws.Cells.Item(31, 13).Value = -387

# Row 32
$ws.Cells.Item(32, 4).Value = 8620
$ws.Cells.Item(32, 5).Value = -219014
$ws.Cells.Item(32, 6).Value = -427113
$ws.Cells.Item(32, 7).Value = -471187
$ws.Cells.Item(32, 8).Value = 295910
$ws.Cells.Item(32, 9).Value = 673893
$ws.Cells.Item(32, 10).Value = 368
$ws.Cells.Item(32, 11).Value = -79042
$ws.Cells.Item(32, 12).Value = -196482
$ws.Cells.Item(32, 13).Value = -1305230

# Row 33
$ws.Cells.Item(33, 4).Value = -116309
$ws.Cells.Item(33, 5).Value = -261761
$ws.Cells.Item(33, 6).Value = 53777
$ws.Cells.Item(33, 7).Value = -1032329
$ws.Cells.Item(33, 8).Value = -148936
$ws.Cells.Item(33, 9).Value = 289902
$ws.Cells.Item(33, 10).Value = -803648
$ws.Cells.Item(33, 11).Value = -1072121
$ws.Cells.Item(33, 12).Value = 931505
$ws.Cells.Item(33, 13).Value = -1328538

# Row 35
$ws.Cells.Item(35, 4).Value = -13
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(35, 9).Value = 0
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 0
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = 0

# Row 36
$ws.Cells.Item(36, 4).Value = "-"
$ws.Cells.Item(36, 5).Value = "-"
$ws.Cells.Item(36, 6).Value = "-"
$ws.Cells.Item(36, 7).Value = "-"
$ws.Cells.Item(36, 8).Value = "-"
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 10).Value = "-"
$ws.Cells.Item(36, 11).Value = "-"
$ws.Cells.Item(36, 12).Value = "-"
$ws.Cells.Item(36, 13).Value = 0

# Row 37
$ws.Cells.Item(37, 4).Value = 0
$ws.Cells.Item(37, 5).Value = 0
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 10).Value = 0
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 12).Value = 0
$ws.Cells.Item(37, 13).Value = 0

# Row 38
$ws.Cells.Item(38, 4).Value = 0
$ws.Cells.Item(38, 5).Value = 0
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 0
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = 0

# Row 39
$ws.Cells.Item(39, 4).Value = 0
$ws.Cells.Item(39, 5).Value = 180000
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 1150000
$ws.Cells.Item(39, 8).Value = 201920
$ws.Cells.Item(39, 9).Value = 797698
$ws.Cells.Item(39, 10).Value = 1129679
$ws.Cells.Item(39, 11).Value = 730605
$ws.Cells.Item(39, 12).Value = 874369
$ws.Cells.Item(39, 13).Value = 2986504

# Row 40
$ws.Cells.Item(40, 4).Value = 70000
$ws.Cells.Item(40, 5).Value = 0
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = -100000
$ws.Cells.Item(40, 9).Value = -319609
$ws.Cells.Item(40, 10).Value = -444733
$ws.Cells.Item(40, 11).Value = 151822
$ws.Cells.Item(40, 12).Value = -594119
$ws.Cells.Item(40, 13).Value = -1503741

# Row 41
$ws.Cells.Item(41, 4).Value = 0
$ws.Cells.Item(41, 5).Value = 0
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 9).Value = -118636
$ws.Cells.Item(41, 10).Value = -50840
$ws.Cells.Item(41, 11).Value = -70902
$ws.Cells.Item(41, 12).Value = -193904
$ws.Cells.Item(41, 13).Value = -348462

# Row 42
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = 0

# Row 43
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(43, 5).Value = 0
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 0
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = 0

# Row 44
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 0
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 9).Value = 0
$ws.Cells.Item(44, 10).Value = 0
$ws.Cells.Item(44, 11).Value = 0
$ws.Cells.Item(44, 12).Value = 0
$ws.Cells.Item(44, 13).Value = 0

# Row 45
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 0
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 0
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = 0
$ws.Cells.Item(45, 13).Value = 0

# Row 46
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 0
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 0
$ws.Cells.Item(46, 9).Value = -360000
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = 0

# Row 47
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = 0

# Row 48
$ws.Cells.Item(48, 4).Value = 0
$ws.Cells.Item(48, 5).Value = 0
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 13).Value = 0

# Row 49
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = 0
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(49, 9).Value = 0
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(49, 12).Value = 0
$ws.Cells.Item(49, 13).Value = 0

# Row 50
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = -14641
$ws.Cells.Item(50, 6).Value = -684
$ws.Cells.Item(50, 7).Value = -247
$ws.Cells.Item(50, 8).Value = -19
$ws.Cells.Item(50, 9).Value = -52266
$ws.Cells.Item(50, 10).Value = -3056
$ws.Cells.Item(50, 11).Value = -695
$ws.Cells.Item(50, 12).Value = -221029
$ws.Cells.Item(50, 13).Value = 53147

# Row 51
$ws.Cells.Item(51, 4).Value = 69987
$ws.Cells.Item(51, 5).Value = 165359
$ws.Cells.Item(51, 6).Value = -684
$ws.Cells.Item(51, 7).Value = 1149753
$ws.Cells.Item(51, 8).Value = 101901
$ws.Cells.Item(51, 9).Value = -52813
$ws.Cells.Item(51, 10).Value = 631050
$ws.Cells.Item(51, 11).Value = 810830
$ws.Cells.Item(51, 12).Value = -134683
$ws.Cells.Item(51, 13).Value = 1187448

# Row 52
$ws.Cells.Item(52, 4).Value = -46322
$ws.Cells.Item(52, 5).Value = -96402
$ws.Cells.Item(52, 6).Value = 53093
$ws.Cells.Item(52, 7).Value = 117424
$ws.Cells.Item(52, 8).Value = -47035
$ws.Cells.Item(52, 9).Value = 237089
$ws.Cells.Item(52, 10).Value = -172598
$ws.Cells.Item(52, 11).Value = -261291
$ws.Cells.Item(52, 12).Value = 796822
$ws.Cells.Item(52, 13).Value = -141090

# Row 53
$ws.Cells.Item(53, 4).Value = 192275
$ws.Cells.Item(53, 5).Value = 144976
$ws.Cells.Item(53, 6).Value = 49576
$ws.Cells.Item(53, 7).Value = 102669
$ws.Cells.Item(53, 8).Value = 230762
$ws.Cells.Item(53, 9).Value = 173058
$ws.Cells.Item(53, 10).Value = 409973
$ws.Cells.Item(53, 11).Value = 237375
$ws.Cells.Item(53, 12).Value = 26534
$ws.Cells.Item(53, 13).Value = 772906

# Row 54
$ws.Cells.Item(54, 4).Value = -977
$ws.Cells.Item(54, 5).Value = 1032
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(54, 7).Value = 10670
$ws.Cells.Item(54, 8).Value = -10670
$ws.Cells.Item(54, 9).Value = -174
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 11).Value = 50450
$ws.Cells.Item(54, 12).Value = -50450
$ws.Cells.Item(54, 13).Value = 76804

# Row 55
$ws.Cells.Item(55, 4).Value = 144976
$ws.Cells.Item(55, 5).Value = 49576
$ws.Cells.Item(55, 6).Value = 102669
$ws.Cells.Item(55, 7).Value = 230762
$ws.Cells.Item(55, 8).Value = 173058
$ws.Cells.Item(55, 9).Value = 409973
$ws.Cells.Item(55, 10).Value = 237375
$ws.Cells.Item(55, 11).Value = 26534
$ws.Cells.Item(55, 12).Value = 772906
$ws.Cells.Item(55, 13).Value = 708620

# Row 56
$ws.Cells.Item(56, 4).Value = 0
$ws.Cells.Item(56, 5).Value = 0
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 11).Value = 0
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 13).Value = 0

